# "write some new for stephen" — duplicate the existing "Subtitle"
# paragraph style definition in word/styles.xml so it appears twice
# (a second, identical <w:style .../> block is inserted immediately
# before the pre-existing one), matching the target OOXML diff.
#
# The Word object model has no supported way to mint a second style
# that shares an existing styleId (Styles.Add just returns the
# existing style when the name/id is already taken), so we reach the
# raw package XML through Document.WordOpenXML, splice in a duplicate
# <w:style> element for "Subtitle", and write it back.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$marker = '<w:style w:type="paragraph" w:styleId="Subtitle">'
$idx = $xml.IndexOf($marker)

if ($idx -ge 0) {
    $closeTag = '</w:style>'
    $endIdx = $xml.IndexOf($closeTag, $idx) + $closeTag.Length
    $subtitleBlock = $xml.Substring($idx, $endIdx - $idx)

    $newXml = $xml.Substring(0, $idx) + $subtitleBlock + $xml.Substring($idx)
    $d.WordOpenXML = $newXml
    Write-Output "Subtitle style duplicated"
} else {
    Write-Output "Subtitle style not found"
}
